$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-17 from 45183 to 45184
$ws.Range("C2:C17").Value = 45184
